$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet: drop the accent in "región" -> "region"
$oldName = $ws.Name
$newName = "des_art_region"
$ws.Name = $newName

# Excel auto-updates most defined names that reference the sheet when it is
# renamed, but make sure every defined name in the workbook is repointed at
# the new sheet name (covers any that weren't automatically rewritten).
foreach ($n in $wb.Names) {
    $ref = $n.RefersTo
    if ($ref -like "*$oldName*") {
        $n.RefersTo = $ref.Replace($oldName, $newName)
    }
}
